$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, pushing existing rows 2-5 down to 3-6
$ws.Rows.Item(2).Insert()

# The insert copies formatting from the row above (the header); clear it so
# the new data row matches the plain (unstyled) look of the other data rows.
$ws.Rows.Item(2).ClearFormats()

# Fill in the new row 2 with the "Ashen Leyndell" boss entry
$ws.Cells.Item(2, 1).Value = "Ashen Leyndell"
$ws.Cells.Item(2, 2).Value = "Major bosses"
$ws.Cells.Item(2, 3).Value = "Dropped by Hoarah Loux. Replaces Remembrance of Hoarah Loux"
